$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Honduras statistics updated and it overtook Burkina Faso in ranking,
# so the two rows swap places (row 102 <-> row 103)
$ws.Range("A102").Value = "Honduras"
$ws.Range("B102").Value = 661
$ws.Range("C102").Value = 34
$ws.Range("D102").Value = 69
$ws.Range("E102").Value = 531
$ws.Range("F102").Value = 10
$ws.Range("G102").Value = 2
$ws.Range("H102").Value = 61

$ws.Range("A103").Value = "Burkina Faso"
$ws.Range("B103").Value = 632
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 453
$ws.Range("E103").Value = 137
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 42

# Camboya row update
$ws.Range("D143").Value = 119
$ws.Range("E143").Value = 3

# Trinidad y Tobago row update
$ws.Range("B144").Value = 116
$ws.Range("C144").Value = 1
$ws.Range("D144").Value = 58
$ws.Range("E144").Value = 50

# Bermudas row update
$ws.Range("E146").Value = 64
$ws.Range("G146").Value = 1
$ws.Range("H146").Value = 6
